$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the descriptive text values in column B to keep naming consistent
# across types, and append a trailing ", " to each (per the shared strings diff).
$ws.Range("B2").Value = "id, providerName(name), tags, refBaseSchema, itemDescription, itemType, "
$ws.Range("B3").Value = "id, tags, refBaseSchema, resourceServer, itemDescription, refDataModel, provider, resourceServerGroup, resourceId, itemType, "
$ws.Range("B4").Value = "id, name, tags, refBaseSchema, itemDescription, resourceServerHTTPAccessURL(uriLink), resourceServerOrg(organizationInfo), coverageRegion, itemType, "
$ws.Range("B5").Value = "id, name, tags, refBaseSchema, resourceServer, itemDescription, refDataModel, provider, itemType, "

# Make the (blue) font used in column A for rows 2-5 bold.
$ws.Range("A2:A5").Font.Bold = $true
